# "todo el pdf andaaaa agregue anexo" — add the Bode-gain (dB) annex column.
#
# Column E = 20*LOG(C/B)  (gain in dB, relative level C over reference B)
# for every data row (1..50), formatted as "0.0" and autofit.
# Also fix the two phase values in D that were left in the wrong units
# (305 -> 305-360 = -55, and 275 -> -85).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D35 / D50 corrections -------------------------------------------------
$ws.Range("D35").Formula = "=305-360"
$ws.Range("D50").Value = -85

# --- New column E: gain in dB ----------------------------------------------
# E1 is entered on its own; E2:E50 filled as one block so Excel records it
# as a single shared-formula group starting at E2 (matches a type-then-
# fill-down authoring flow).
$ws.Range("E1").Formula = "=20*LOG(C1/B1)"
$ws.Range("E2:E50").Formula = "=20*LOG(C2/B2)"

$ws.Range("E1:E50").NumberFormat = "0.0"

# Size column E to fit its contents (best-fit width for "0.0"-formatted
# gain values, e.g. "-22.6").
$ws.Columns("E:E").ColumnWidth = 11.45

# --- Selection mirrors the final used range ---------------------------------
$ws.Range("A1:E50").Select()
